# Update "想去人数" (attendance) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 298
    $ws.Range("F3").Value = 94
    $ws.Range("F4").Value = 1169
}
